$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training metrics for epochs 15-21 (replacing previous epochs 1-12)
$data = @(
    @(15, 0.83,  0.84375,            0.8100000000000001, 0.826530612244898,  25.28125,  1.46875, 7644233728,  260.9365339279175),
    @(16, 0.825, 0.8350515463917526, 0.8100000000000001, 0.8223350253807107, 22.21875,  2.59375, 682217472,   332.3033769130707),
    @(17, 0.83,  0.84375,            0.8100000000000001, 0.826530612244898,  22.578125, 5.640625, -2554433536, 392.7677466869354),
    @(18, 0.835, 0.845360824742268,  0.82,                0.83248730964467,   21.28125,  4.96875, 724201472,   369.5096917152405),
    @(19, 0.835, 0.845360824742268,  0.82,                0.83248730964467,   20.703125, 7.078125, -773591040, 407.268709897995),
    @(20, 0.835, 0.845360824742268,  0.82,                0.83248730964467,   20.25,     5.53125, -95969280,   398.269159078598),
    @(21, 0.835, 0.845360824742268,  0.82,                0.83248730964467,   20.3125,   7.140625, 195874816,  421.8510708808899)
)

# Remove the trailing rows (old epochs 8-12), since the new data only has 7 rows (was 12)
$ws.Range("A9:I13").EntireRow.Delete() | Out-Null

# Write new values into rows 2-8, columns A-I
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
